$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 497.5207405739714
$ws.Range("C2").Value = 11.92123037940973
$ws.Range("D2").Value = 231.5049404543167
